# Generate Report for Archive
#
# The three rows for b7445b20..., 36e46dc5... and a75ec3ab... get re-ordered:
#   a75ec3ab...   moves up to row 4
#   b7445b20...   moves down to row 5
#   36e46dc5...   moves down to row 6
#
# Cell values move with the rows, but (matching observed Excel hyperlink
# behaviour) each hyperlink's target Address stays anchored to its
# worksheet position -- only the displayed text is refreshed to track the
# new cell content.

$wb = $excel.ActiveWorkbook

function Update-Row {
    param($ws, $rowIndex, $values, $hyperlinkCols)

    # values is a hashtable col-letter -> new text
    foreach ($col in $values.Keys) {
        $addr = ($col + $rowIndex)
        $ws.Range($addr).Value = $values[$col]
    }

    # Collect this sheet's hyperlinks once (Hyperlinks.Item(i) is unreliable
    # in this host, so capture via foreach and index into the array).
    $links = @()
    foreach ($h in $ws.Hyperlinks) {
        $links += $h
    }

    foreach ($col in $hyperlinkCols) {
        $targetAddr = ($col + $rowIndex)
        foreach ($h in $links) {
            if ($h.Range.Address($false, $false) -eq $targetAddr) {
                $h.TextToDisplay = $values[$col]
            }
        }
    }
}

# ---- Sheet "Overview" ----
$ws1 = $wb.Worksheets.Item("Overview")

Update-Row $ws1 4 @{ "A" = "a75ec3ab-0cdb-478b-9715-f84e05f42278.md"; "B" = "In Translation"; "C" = "In Translation"; "D" = "2016-32-19 16:32:08" } @("A")
Update-Row $ws1 5 @{ "A" = "b7445b20-9b0c-4b5a-bc9a-d6eb9e5c884b.md"; "B" = "In Translation"; "C" = "In Translation"; "D" = "2016-31-19 16:31:30" } @("A")
Update-Row $ws1 6 @{ "A" = "36e46dc5-c445-48ce-af7f-65d1a668a69e.md"; "B" = "Ready for handoff"; "C" = "Ready for handoff"; "D" = "2016-32-19 16:32:27" } @("A")

# ---- Sheet "zh-cn" ----
$ws2 = $wb.Worksheets.Item("zh-cn")

Update-Row $ws2 4 @{ "A" = "a75ec3ab-0cdb-478b-9715-f84e05f42278.md"; "C" = "In Translation"; "D" = "a75ec3ab-0cdb-478b-9715-f84e05f42278.9c5fbd2db760637f650df4b0b3d3f7ba3f6131f7.zh-cn.xlf"; "E" = "2016-03-19 16:32:05" } @("A", "D")
Update-Row $ws2 5 @{ "A" = "b7445b20-9b0c-4b5a-bc9a-d6eb9e5c884b.md"; "C" = "In Translation"; "D" = "b7445b20-9b0c-4b5a-bc9a-d6eb9e5c884b.878e7a760d652024fb3855f3435cf0bb2bd57238.zh-cn.xlf"; "E" = "2016-03-19 16:31:27" } @("A", "D")
Update-Row $ws2 6 @{ "A" = "36e46dc5-c445-48ce-af7f-65d1a668a69e.md"; "C" = "Ready for handoff"; "D" = "36e46dc5-c445-48ce-af7f-65d1a668a69e.11c27fc6f0a38c82586726e8f9eac026f0b9672b.zh-cn.xlf"; "E" = "2016-03-19 16:32:25" } @("A", "D")

# ---- Sheet "de-de" ----
$ws3 = $wb.Worksheets.Item("de-de")

Update-Row $ws3 4 @{ "A" = "a75ec3ab-0cdb-478b-9715-f84e05f42278.md"; "C" = "In Translation"; "D" = "a75ec3ab-0cdb-478b-9715-f84e05f42278.9c5fbd2db760637f650df4b0b3d3f7ba3f6131f7.de-de.xlf"; "E" = "2016-03-19 16:32:08" } @("A", "D")
Update-Row $ws3 5 @{ "A" = "b7445b20-9b0c-4b5a-bc9a-d6eb9e5c884b.md"; "C" = "In Translation"; "D" = "b7445b20-9b0c-4b5a-bc9a-d6eb9e5c884b.878e7a760d652024fb3855f3435cf0bb2bd57238.de-de.xlf"; "E" = "2016-03-19 16:31:30" } @("A", "D")
Update-Row $ws3 6 @{ "A" = "36e46dc5-c445-48ce-af7f-65d1a668a69e.md"; "C" = "Ready for handoff"; "D" = "36e46dc5-c445-48ce-af7f-65d1a668a69e.11c27fc6f0a38c82586726e8f9eac026f0b9672b.de-de.xlf"; "E" = "2016-03-19 16:32:27" } @("A", "D")
